# Update the "repaymentstrategy" value on the ProductLoanInput sheet
# from "RBI (India)" to "Overdue/Due Fee/Int,Principal"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Move/update the active selection to B17 to match the saved view state
$ws.Activate()
$ws.Range("B17").Select()
